$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (TSM)
$ws.Range("D2").Value = 293.25
$ws.Range("E2").Value = 60.7
$ws.Range("F2").Value = 1.13
$ws.Range("I2").Value = 63
$ws.Range("K2").Value = 58.9
$ws.Range("N2").Value = 52.47848103381103

# Row 3 (ASML)
$ws.Range("N3").Value = 52.47848103381103

# Row 4 (QCOM)
$ws.Range("D4").Value = 174.49
$ws.Range("F4").Value = 5.66
$ws.Range("N4").Value = 52.47848103381103

# Row 5 (NVDA)
$ws.Range("D5").Value = 184.34
$ws.Range("E5").Value = 47.5
$ws.Range("F5").Value = 2.26
$ws.Range("N5").Value = 52.47848103381103

# Row 6 (AMD)
$ws.Range("D6").Value = 218.94
$ws.Range("F6").Value = 2.19
$ws.Range("N6").Value = 52.47848103381103
